# Append: 2025-12-01 02:09 JST
# Update the "取得日時" (acquisition timestamp) column (A) for all existing
# data rows on the "ランサーズ" sheet from 2025-12-01 01:44:37 to 2025-12-01 02:09:15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-01 02:09:15"

for ($row = 2; $row -le 17; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
